$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Задачи")

# TASK_14 (row 15): fill in completion date (column E = "заврешенно")
$ws.Range("E15").Value = Get-Date -Year 2014 -Month 12 -Day 12 -Hour 17 -Minute 8 -Second 0

# TASK_15 (row 16): fill in start date (D) and completion date (E)
$ws.Range("D16").Value = Get-Date -Year 2014 -Month 12 -Day 14 -Hour 12 -Minute 35 -Second 0
$ws.Range("E16").Value = Get-Date -Year 2014 -Month 12 -Day 14 -Hour 14 -Minute 47 -Second 0

# Update selection to D17
$ws.Range("D17").Select()
